# 10月学習計画書 - update study plan rows for Node.js / React progress
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("学習計画書")

# Row 6: narrow the first Node.js chunk from "1-15" to "1-5"
$ws.Range("F6").Value = "Nodejs MongoDb Express 零基础 入门 实战 视频教程 1-5・練習"

# Row 7: next Node.js chunk "16-30" -> "6-15"; mark complete; add remark
$ws.Range("F7").Value = "Nodejs MongoDb Express 零基础 入门 实战 视频教程 6-15・練習"
$ws.Range("G7").NumberFormat = "0%"
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = "【尚硅谷】NodeJS全套教程　11-15"

# Row 8: next Node.js chunk "31-39" -> "16-23"; mark complete
$ws.Range("F8").Value = "Nodejs MongoDb Express 零基础 入门 实战 视频教程 16-23・練習"
$ws.Range("G8").NumberFormat = "0%"
$ws.Range("G8").Value = 1

# Row 9 (Saturday) now has a study session too
$ws.Range("C9").Value = 0.39583333333333331
$ws.Range("D9").Value = 0.77083333333333337
$ws.Range("E9").Value = 0.041666666666666664
$ws.Range("F9").Value = "Nodejs MongoDb Express 零基础 入门 实战 视频教程 24-32・練習"
$ws.Range("H9").Value = "【尚硅谷】NodeJS全套教程　16-20"

# Row 10 (Sunday) now has a study session too
$ws.Range("C10").Value = 0.39583333333333331
$ws.Range("D10").Value = 0.77083333333333337
$ws.Range("E10").Value = 0.041666666666666664
$ws.Range("F10").Value = "Nodejs MongoDb Express 零基础 入门 实战 视频教程 33-39・練習"

# Move the active selection to F8, matching the saved cursor position
$ws.Range("F8").Select()
